# Update Peak_Area_SO2 (column C) values for rows 2-9 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 556.8820601375211
$ws.Range("C3").Value = 45.47286828470842
$ws.Range("C4").Value = 180.0007210602771
$ws.Range("C5").Value = 70.08344443838996
$ws.Range("C6").Value = 1193.886885480184
$ws.Range("C7").Value = 336.8006864955747
$ws.Range("C8").Value = 361.1135807444104
$ws.Range("C9").Value = 88.45279340605903
